# feat: add 2022-Q4 data
#
# Before:  总计 | 2022-Q3 | 2022-Q2
# After:   总计 | 2022-Q4 | 2022-Q3 | 2022-Q2
#
# A new "2022-Q4" tab is inserted right after "总计" (reusing the
# worksheet that used to hold "2022-Q3"'s numbers), and the former
# "2022-Q3" numbers are preserved on a freshly duplicated tab so the
# old quarter keeps its data. The "总计" summary sheet gets a new row
# for 2022-Q4, pushing the Q3/Q2 rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Duplicate the "2022-Q3" sheet BEFORE touching its numbers, so the
#    copy keeps the old Q3 figures. The copy is placed immediately
#    after the source sheet.
# ---------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Item(2)           # currently "2022-Q3" (old data)
$wsQ3.Copy([System.Reflection.Missing]::Value, $wsQ3)

$wsQ3Copy = $wb.Worksheets.Item(3)       # the duplicate, named "2022-Q3 (2)"

# ---------------------------------------------------------------
# 2) Repurpose the original sheet as "2022-Q4" (rename it out of the way
#    first so the duplicate can reclaim the "2022-Q3" name) and refresh
#    its figures.
# ---------------------------------------------------------------
$wsQ3.Name = "2022-Q4"
$wsQ3Copy.Name = "2022-Q3"               # this tab keeps the OLD Q3 numbers

$rows = @(2, 3, 4)
$newD = "12.08"
$newE = "89.83"
$newF = "4.11"
$newG = "0.4965"

# Keep these text-looking numbers stored as text (matches the source
# file, where D:G are inline strings, not numeric cells).
$wsQ3.Range("D2:G4").NumberFormat = "@"

foreach ($r in $rows) {
    $wsQ3.Cells.Item($r, 4).Value = $newD
    $wsQ3.Cells.Item($r, 5).Value = $newE
    $wsQ3.Cells.Item($r, 6).Value = $newF
    $wsQ3.Cells.Item($r, 7).Value = $newG
}

# ---------------------------------------------------------------
# 3) Update the "总计" (summary) sheet: insert a 2022-Q4 row, pushing
#    the existing Q3/Q2 rows down by one.
# ---------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

# Push row 3 (2022-Q2) down into row 4, copying formatting from row 3.
$wsTotal.Cells.Item(3, 1).Copy()
$wsTotal.Cells.Item(4, 1).PasteSpecial(-4122)   # xlPasteFormats
$wsTotal.Cells.Item(4, 1).Value = 2
$wsTotal.Cells.Item(4, 2).Value = $wsTotal.Cells.Item(3, 2).Value()
$wsTotal.Cells.Item(4, 3).Value = $wsTotal.Cells.Item(3, 3).Value()
$wsTotal.Cells.Item(4, 4).Value = $wsTotal.Cells.Item(3, 4).Value()

# Push row 2 (2022-Q3) down into row 3 (A3 already has the right style).
$wsTotal.Cells.Item(3, 1).Value = 1
$wsTotal.Cells.Item(3, 2).Value = $wsTotal.Cells.Item(2, 2).Value()
$wsTotal.Cells.Item(3, 3).Value = $wsTotal.Cells.Item(2, 3).Value()
$wsTotal.Cells.Item(3, 4).Value = $wsTotal.Cells.Item(2, 4).Value()

# Write the new 2022-Q4 row into row 2 (A2 keeps its existing style/value).
$wsTotal.Cells.Item(2, 2).Value = "2022-Q4"
$wsTotal.Cells.Item(2, 3).Value = 3
$wsTotal.Cells.Item(2, 4).Value = 1.49
